# Auto-generated edit script applying the Sargatanas_Profits.xlsx commit diff
# (scheduled-runner profit recompute touching ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 28
$ws.Range("H28").Value = 1051.7059
$ws.Range("J28").Value = 610.2222
$ws.Range("L28").Value = 610.2222
$ws.Range("N28").Value = -1580.2222
# row 32
$ws.Range("H32").Value = 3681.3845
$ws.Range("I32").Value = 2917.6
$ws.Range("J32").Value = 4158.75
$ws.Range("K32").Value = 2917.6
$ws.Range("L32").Value = 4158.75
$ws.Range("M32").Value = -2591.6
$ws.Range("N32").Value = -4810.75
# row 62
$ws.Range("H62").Value = 66720268
$ws.Range("I62").Value = 200002400
$ws.Range("J62").Value = 79202.5
$ws.Range("K62").Value = 200002400
$ws.Range("L62").Value = 79202.5
$ws.Range("M62").Value = -200001776
$ws.Range("N62").Value = -80450.5
# row 65
$ws.Range("H65").Value = 66720268
$ws.Range("I65").Value = 200002400
$ws.Range("J65").Value = 79202.5
$ws.Range("K65").Value = 1000012000
$ws.Range("L65").Value = 396012.5
$ws.Range("M65").Value = -1000008880
$ws.Range("N65").Value = -402252.5
# row 70
$ws.Range("H70").Value = 72920370
$ws.Range("I70").Value = 83336460
$ws.Range("K70").Value = 250009380
$ws.Range("M70").Value = -250009110
# row 73
$ws.Range("H73").Value = 72920370
$ws.Range("I73").Value = 83336460
$ws.Range("K73").Value = 250009380
$ws.Range("M73").Value = -250008444
# row 100
$ws.Range("H100").Value = 2454.7856
$ws.Range("I100").Value = 1840.8
$ws.Range("K100").Value = 1840.8
$ws.Range("M100").Value = -1299.8
# row 103
$ws.Range("H103").Value = 921.26666
$ws.Range("I103").Value = 559.6667
$ws.Range("J103").Value = 1011.6667
$ws.Range("K103").Value = 1679.0001
$ws.Range("L103").Value = 3035.0001
$ws.Range("M103").Value = -1093.0001
$ws.Range("N103").Value = -4207.0001
# row 107
$ws.Range("H107").Value = 11030500
$ws.Range("J107").Value = 35716210
$ws.Range("L107").Value = 35716210
$ws.Range("N107").Value = -35720050
# row 125
$ws.Range("H125").Value = 83339590
$ws.Range("J125").Value = 8383.75
$ws.Range("L125").Value = 75453.75
$ws.Range("N125").Value = -80373.75
# row 135
$ws.Range("H135").Value = 476783.25
$ws.Range("J135").Value = 348.25
$ws.Range("L135").Value = 3134.25
$ws.Range("N135").Value = -8204.25
# row 137
$ws.Range("H137").Value = 2417.0881
$ws.Range("I137").Value = 1938.9
$ws.Range("K137").Value = 5816.700000000001
$ws.Range("M137").Value = -3266.700000000001

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 1457056.8
$ws.Range("I32").Value = 1565674.8
$ws.Range("K32").Value = 1565674.8
$ws.Range("M32").Value = -1565387.8
# row 44
$ws.Range("H44").Value = 12044
$ws.Range("I44").Value = 12044
$ws.Range("K44").Value = 12044
$ws.Range("M44").Value = -11556
# row 52
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = ""
# row 61
$ws.Range("H61").Value = 5065.6777
$ws.Range("I61").Value = 2965.8206
$ws.Range("J61").Value = 9160.4
$ws.Range("K61").Value = 2965.8206
$ws.Range("L61").Value = 9160.4
$ws.Range("M61").Value = -2753.8206
$ws.Range("N61").Value = -9584.4
# row 110
$ws.Range("H110").Value = 20834752
$ws.Range("I110").Value = 1385.2727
$ws.Range("K110").Value = 1385.2727
$ws.Range("M110").Value = 659.7273
# row 122
$ws.Range("H122").Value = 2609.7083
$ws.Range("I122").Value = 1191.2106
$ws.Range("K122").Value = 3573.6318
$ws.Range("M122").Value = -1123.6318
# row 136
$ws.Range("H136").Value = 5065.6777
$ws.Range("I136").Value = 2965.8206
$ws.Range("J136").Value = 9160.4
$ws.Range("K136").Value = 8897.461800000001
$ws.Range("L136").Value = 27481.2
$ws.Range("M136").Value = -6347.461800000001
$ws.Range("N136").Value = -32581.2

$ws = $wb.Worksheets.Item("BSM")
# row 99
$ws.Range("H99").Value = 2275515.8
$ws.Range("I99").Value = 2566.1562
$ws.Range("J99").Value = 11367314
$ws.Range("K99").Value = 2566.1562
$ws.Range("L99").Value = 11367314
$ws.Range("M99").Value = -1068.1562
$ws.Range("N99").Value = -11370310

$ws = $wb.Worksheets.Item("CRP")
# row 68
$ws.Range("H68").Value = 85000
$ws.Range("J68").Value = 85000
$ws.Range("L68").Value = 85000
$ws.Range("N68").Value = -86498
# row 71
$ws.Range("H71").Value = 85000
$ws.Range("J71").Value = 85000
$ws.Range("L71").Value = 255000
$ws.Range("N71").Value = -262488
# row 99
$ws.Range("H99").Value = 4521.737
$ws.Range("I99").Value = 2849.9285
$ws.Range("K99").Value = 2849.9285
$ws.Range("M99").Value = -1351.9285
# row 107
$ws.Range("H107").Value = 1428.625
$ws.Range("I107").Value = 1294.5
$ws.Range("J107").Value = 1562.75
$ws.Range("K107").Value = 1294.5
$ws.Range("L107").Value = 1562.75
$ws.Range("M107").Value = 625.5
$ws.Range("N107").Value = -5402.75
# row 126
$ws.Range("H126").Value = 4521.737
$ws.Range("I126").Value = 2849.9285
$ws.Range("K126").Value = 8549.7855
$ws.Range("M126").Value = -6079.7855
# row 134
$ws.Range("H134").Value = 4610.8076
$ws.Range("I134").Value = 1309.5358
$ws.Range("J134").Value = 8462.291999999999
$ws.Range("K134").Value = 3928.6074
$ws.Range("L134").Value = 25386.876
$ws.Range("M134").Value = -1393.6074
$ws.Range("N134").Value = -30456.876

$ws = $wb.Worksheets.Item("CUL")
# row 86
$ws.Range("H86").Value = 417.72223
$ws.Range("I86").Value = 408
$ws.Range("J86").Value = 466.33334
$ws.Range("K86").Value = 1224
$ws.Range("L86").Value = 1399.00002
$ws.Range("M86").Value = -38
$ws.Range("N86").Value = -3771.00002
# row 89
$ws.Range("H89").Value = 417.72223
$ws.Range("I89").Value = 408
$ws.Range("J89").Value = 466.33334
$ws.Range("K89").Value = 3672
$ws.Range("L89").Value = 4197.00006
$ws.Range("M89").Value = 2256
$ws.Range("N89").Value = -16053.00006

$ws = $wb.Worksheets.Item("GSM")
# row 102
$ws.Range("H102").Value = 4185.8184
$ws.Range("I102").Value = 4427.6
$ws.Range("K102").Value = 4427.6
$ws.Range("M102").Value = -2805.6
# row 117
$ws.Range("H117").Value = 41309.2
$ws.Range("J117").Value = 41309.2
$ws.Range("L117").Value = 41309.2
$ws.Range("N117").Value = -48193.2
# row 126
$ws.Range("H126").Value = 5295.6665
$ws.Range("I126").Value = 5927.857
$ws.Range("J126").Value = 4410.6
$ws.Range("K126").Value = 17783.571
$ws.Range("L126").Value = 13231.8
$ws.Range("M126").Value = -15313.571
$ws.Range("N126").Value = -18171.8

$ws = $wb.Worksheets.Item("LTW")
# row 10
$ws.Range("H10").Value = 266.66666
# row 16
$ws.Range("H16").Value = 1134.4286
$ws.Range("I16").Value = 788.2
$ws.Range("K16").Value = 788.2
$ws.Range("M16").Value = -618.2
# row 22
$ws.Range("H22").Value = 2331.1667
$ws.Range("I22").Value = 358.66666
$ws.Range("K22").Value = 358.66666
$ws.Range("M22").Value = -63.66665999999998
# row 27
$ws.Range("H27").Value = 2331.1667
$ws.Range("I27").Value = 358.66666
$ws.Range("K27").Value = 358.66666
$ws.Range("M27").Value = -251.66666
# row 40
$ws.Range("H40").Value = 6158.5884
$ws.Range("I40").Value = 4948.5
$ws.Range("K40").Value = 4948.5
$ws.Range("M40").Value = -4812.5
# row 93
$ws.Range("H93").Value = 5070.7617
$ws.Range("I93").Value = 6162.75
$ws.Range("K93").Value = 6162.75
$ws.Range("M93").Value = -4914.75
# row 100
$ws.Range("H100").Value = 3870
$ws.Range("I100").Value = 2741.5
$ws.Range("K100").Value = 2741.5
$ws.Range("M100").Value = -2200.5
# row 132
$ws.Range("H132").Value = 11117745
$ws.Range("I132").Value = 27780308
$ws.Range("J132").Value = 9370.296
$ws.Range("K132").Value = 83340924
$ws.Range("L132").Value = 28110.888
$ws.Range("M132").Value = -83338394
$ws.Range("N132").Value = -33170.888
# row 136
$ws.Range("H136").Value = 7905.2974
$ws.Range("I136").Value = 1779.8667
$ws.Range("J136").Value = 12081.728
$ws.Range("K136").Value = 5339.6001
$ws.Range("L136").Value = 36245.18399999999
$ws.Range("M136").Value = -2789.6001
$ws.Range("N136").Value = -41345.18399999999

$ws = $wb.Worksheets.Item("WVR")
# row 59
$ws.Range("H59").Value = 50000
$ws.Range("J59").Value = 50000
$ws.Range("L59").Value = 50000
$ws.Range("M59").Value = -51476
# row 62
$ws.Range("H62").Value = 2999.5
$ws.Range("J62").Value = 2999
$ws.Range("L62").Value = 2999
$ws.Range("N62").Value = -4247
# row 65
$ws.Range("H65").Value = 2999.5
$ws.Range("J65").Value = 2999
$ws.Range("L65").Value = 14995
$ws.Range("N65").Value = -21235
# row 126
$ws.Range("H126").Value = 2751.0334
$ws.Range("I126").Value = 1686.4
$ws.Range("J126").Value = 3815.6667
$ws.Range("K126").Value = 5059.200000000001
$ws.Range("L126").Value = 11447.0001
$ws.Range("M126").Value = -2589.200000000001
$ws.Range("N126").Value = -16387.0001

